# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Apio" (Vega Modelo de Temuco)
# right above the existing row 170, shifting all subsequent rows down by
# one (old row 170 -> 171, ..., old row 194 -> 195).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 170; Excel pushes rows 170:194 down to 171:195.
$ws.Rows(170).Insert()

# Populate the newly inserted row 170 with the new weekly record.
$ws.Range("A170").Value = 10
$ws.Range("B170").Value = "Vega Modelo de Temuco"
$ws.Range("C170").Value = "La Araucanía"
$ws.Range("D170").Value = 44491
$ws.Range("E170").Value = 9
$ws.Range("F170").Value = 100112017
$ws.Range("G170").Value = "Apio"
$ws.Range("H170").Value = "Americana (o)"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 145
$ws.Range("K170").Value = 8000
$ws.Range("L170").Value = 9000
$ws.Range("M170").Value = 8448
$ws.Range("N170").Value = "$/docena de matas"
$ws.Range("O170").Value = "Provincia del Elquí"
$ws.Range("P170").Value = 1408
$ws.Range("Q170").Value = 6
$ws.Range("R170").Value = "Hortaliza"
